$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every employee's email address: old "name@surname.it" style
# addresses are swapped for a uniform "name@gmail.com" address. Values are
# written in row order (skipping the two rows that already carry
# hyperlinks, E14 and E19, until the end) so the shared-string table is
# rebuilt in the same order the original authoring tool produced.
$ws.Range("E2").Value = "giovanni@gmail.com"
$ws.Range("E3").Value = "laura@gmail.com"
$ws.Range("E4").Value = "alessandro@gmail.com"
$ws.Range("E5").Value = "simona@gmail.com"
$ws.Range("E6").Value = "daniele@gmail.com"
$ws.Range("E7").Value = "elena@gmail.com"
$ws.Range("E8").Value = "roberto@gmail.com"
$ws.Range("E9").Value = "silvia@gmail.com"
$ws.Range("E10").Value = "marco@gmail.com"
$ws.Range("E11").Value = "valentina@gmail.com"
$ws.Range("E12").Value = "nicola@gmail.com"
$ws.Range("E13").Value = "giorgia@gmail.com"
$ws.Range("E15").Value = "elisa@gmail.com"
$ws.Range("E16").Value = "francesco@gmail.com"
$ws.Range("E17").Value = "martina@gmail.com"
$ws.Range("E18").Value = "paolo@gmail.com"
$ws.Range("E20").Value = "riccardo@gmail.com"
$ws.Range("E21").Value = "riccardino@gmail.com"
$ws.Range("E14").Value = "luigi@gmail.com"
$ws.Range("E19").Value = "federica@gmail.com"

# Hyperlinks move: Riccardo Gatti (E20) no longer has a mailto link, while
# Luigi Fontana (E14) and Federica Caruso (E19) gain one. The underlying
# engine's Hyperlinks.Delete() clears the whole sheet collection, so drop
# everything once and rebuild the three links that should remain/appear,
# in the same order as the saved workbook (E21, E14, E19).
$ws.Range("E20").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E21"), "mailto:riccardino@gmail.com")
$ws.Range("E21").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E14"), "mailto:luigi@gmail.com")
$ws.Range("E14").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E19"), "mailto:federica@gmail.com")
$ws.Range("E19").Style = "Hyperlink"

# Restore the selection cursor position recorded in the saved file.
[void]$ws.Range("E20").Select()
